$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D6").Value = 0.005
$ws.Range("E2:E6").Value = 0.01580442584641887
